# Add "2022-Q3" data: new sheet inserted right after "总计" (before "2021-Q2"),
# and a new row on the "总计" summary sheet for the new quarter.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal TEXT (no leading-quote artifact, no
# leftover number-format on the cell) by temporarily forcing a text
# number-format, assigning the value, then clearing the format again.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before "2021-Q2".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2021-Q2"))
$newSheet.Name = "2022-Q3"

# NOTE: worksheet references resolve by tab position, so anything
# fetched *before* the Add() above can now point at the wrong tab.
# Re-resolve "2021-Q2" by name now that the sheet collection has
# changed, and use it (not a stale handle) as the copy source.
$q2Sheet = $wb.Worksheets.Item("2021-Q2")

# Copy header row + first data row (with their styles) from "2021-Q2"
# so the new sheet's formatting (bold/centered/bordered header, styled
# column A) matches the existing quarter sheets exactly.
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1"))
$q2Sheet.Range("A2:H2").Copy($newSheet.Range("A2"))

# Fill in the 2022-Q3 fund-holding data.
$newSheet.Range("D1").Value = "基金规模"

$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "159628"
Set-TextValue $newSheet.Range("C2") "万家国证2000ETF"
Set-TextValue $newSheet.Range("D2") "2.90"
Set-TextValue $newSheet.Range("E2") "97.72"
Set-TextValue $newSheet.Range("F2") "0.44"
Set-TextValue $newSheet.Range("G2") "0.0128"
$newSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Add the matching row to the "总计" (total) sheet, keeping the
#    existing 2021-Q2 / 2021-Q1 rows intact (just shifted down).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row should carry no special formatting except column A,
# which (like the rows below it) uses the bold/centered style.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# 3. Restore "2021-Q1" as the active/selected tab (unchanged from the
#    original workbook).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
